# Update cryptocurrency price/volume snapshot values (Price column D, Volume(1h) column E)
# Values are kept as text (leading apostrophe forces text entry) to match the
# original inlineStr string cell type used by the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.77"
$ws.Range("E2").Value = "'-4.88%"
$ws.Range("D3").Value = "'48.93"
$ws.Range("E3").Value = "'-1.76%"
$ws.Range("D4").Value = "'5.175"
$ws.Range("E4").Value = "'-3.34%"
$ws.Range("D5").Value = "'0.07734"
$ws.Range("E5").Value = "'-5.30%"
$ws.Range("D6").Value = "'4.503"
$ws.Range("E6").Value = "'-2.33%"
$ws.Range("E7").Value = "'13.82%"
$ws.Range("D8").Value = "'1.547"
$ws.Range("E8").Value = "'-7.23%"
$ws.Range("D9").Value = "'0.1220"
$ws.Range("E9").Value = "'-9.89%"
$ws.Range("D10").Value = "'0.1933"
$ws.Range("E10").Value = "'-1.12%"
$ws.Range("D11").Value = "'0.04672"
$ws.Range("E11").Value = "'2.50%"
$ws.Range("D12").Value = "'0.09273"
$ws.Range("E12").Value = "'-3.32%"
$ws.Range("D13").Value = "'0.1046"
$ws.Range("E13").Value = "'-0.30%"
$ws.Range("D14").Value = "'0.001267"
$ws.Range("E14").Value = "'-3.90%"
$ws.Range("D15").Value = "'0.04181"
$ws.Range("E15").Value = "'-3.00%"
$ws.Range("D16").Value = "'0.005813"
$ws.Range("E16").Value = "'-2.24%"
$ws.Range("D17").Value = "'3.327"
$ws.Range("E17").Value = "'-2.10%"
$ws.Range("E18").Value = "'-8.02%"
$ws.Range("D19").Value = "'0.3487"
$ws.Range("D20").Value = "'8.018"
$ws.Range("E20").Value = "'-2.16%"
$ws.Range("E21").Value = "'-5.63%"
$ws.Range("E22").Value = "'-0.40%"
$ws.Range("D23").Value = "'0.001277"
$ws.Range("E23").Value = "'-2.02%"
$ws.Range("D24").Value = "'0.004083"
$ws.Range("E24").Value = "'-4.35%"
$ws.Range("E25").Value = "'0.27%"
$ws.Range("E26").Value = "'-4.05%"
$ws.Range("E38").Value = "'-7.56%"
$ws.Range("D39").Value = "'0.05796"
$ws.Range("E39").Value = "'4.25%"
$ws.Range("D40").Value = "'0.01076"
$ws.Range("E40").Value = "'70.74%"
$ws.Range("D41").Value = "'0.007958"
$ws.Range("E41").Value = "'3.41%"
$ws.Range("E42").Value = "'-2.19%"
$ws.Range("D43").Value = "'0.008326"
$ws.Range("E43").Value = "'8.33%"
$ws.Range("D44").Value = "'0.007693"
$ws.Range("E44").Value = "'-4.56%"
$ws.Range("D45").Value = "'0.3359"
$ws.Range("E45").Value = "'-4.44%"
$ws.Range("D46").Value = "'0.00007002"
$ws.Range("E46").Value = "'3.41%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.19%"
$ws.Range("D48").Value = "'0.05670"
$ws.Range("E48").Value = "'-7.53%"
$ws.Range("E49").Value = "'0.14%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.19%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.19%"
